$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years = 2012..2021
for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $years[$i]
}
